# Addition of filtering option + fixing orientation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "revenue" header and values in column E
$ws.Range("E1").Value = "revenue"

$revenue = @(100, 125, 150, 900, 700, 1200, 111, 90, 400, 600, 800, 888, 123, 145, 1556, 1720, 677, 870, 875, 345, 90, 1200)

for ($i = 0; $i -lt $revenue.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $revenue[$i]
}

# Fix selection/orientation: active cell now E7
$ws.Range("E7").Select()
